$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.884.27'
$ws.Range("E2").Value = '  -1.55%  '

$ws.Range("D3").Value = '3.870.01'
$ws.Range("E3").Value = '  -1.62%  '

$ws.Range("E4").Value = '  +0.07%  '

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '598.58'
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = '  -0.95%  '

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.82'
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = '  +1.41%  '

$ws.Range("D7").Value = '3.868.60'
$ws.Range("E7").Value = '  -1.63%  '

$ws.Range("E8").Value = '  +0.11%  '

$ws.Range("E9").Value = '  -0.61%  '

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = '  -3.51%  '

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.38'
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = '  -1.91%  '

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.456'
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = '  -2.33%  '

$ws.Range("E13").Value = '  +0.85%  '

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.90'
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("D15").Value = '4.520.75'
$ws.Range("E15").Value = '  -1.49%  '

$ws.Range("D16").Value = '3.872.22'
$ws.Range("E16").Value = '  -0.94%  '

$ws.Range("D17").Value = '68.052.04'
$ws.Range("E17").Value = '  -1.33%  '

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.09'
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = '  +3.97%  '

$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.33'
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = '  -2.42%  '

$ws.Range("E20").Value = '  -0.45%  '

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.84'
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = '  -1.66%  '

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '466.31'
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = '  -5.69%  '

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.737'
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = '  +0.78%  '

$ws.Range("E24").Value = '  -3.99%  '

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.16'
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = '  -2.06%  '

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.23'
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = '  -1.86%  '

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.05'
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = '  -1.12%  '

$ws.Range("E28").Value = '  -0.09%  '

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.95'
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = '  -2.47%  '

$ws.Range("E30").Value = '  -1.10%  '

$ws.Range("D31").Value = '4.023.10'
$ws.Range("E31").Value = '  -1.48%  '

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.69'
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = '  -1.75%  '

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.30'
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = '  -3.98%  '

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.16'
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = '  -3.10%  '

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.42'
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = '  -1.17%  '

$ws.Range("D36").Value = '3.841.31'
$ws.Range("E36").Value = '  -1.31%  '

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.74'
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = '  +13.53%  '

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = '  -2.98%  '

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = '  -2.15%  '

$ws.Range("E40").Value = '  +0.11%  '

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.89'
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = '  -2.03%  '

$ws.Range("E42").Value = '  +0.10%  '

$ws.Range("E43").Value = '  -2.92%  '

$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000300'
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = '  +10.72%  '

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '423.17'
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = '  -3.89%  '

$ws.Range("E46").Value = '  -1.20%  '

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.61'
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = '  +0.00%  '

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '47.10'
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = '  -2.18%  '

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '27.39'
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = '  +6.15%  '

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '143.27'
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = '  +0.13%  '
